$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("C2:C21").Value = "MS"
$ws.Range("C13").Select() | Out-Null
